$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Cells.Item(76, 8).Value = 35334
$ws.Cells.Item(76, 9).Value = 35334
$ws.Cells.Item(76, 11).Value = 35334
$ws.Cells.Item(76, 13).Value = -35019
# Row 79
$ws.Cells.Item(79, 8).Value = 35334
$ws.Cells.Item(79, 9).Value = 35334
$ws.Cells.Item(79, 11).Value = 35334
$ws.Cells.Item(79, 13).Value = -34242
# Row 87
$ws.Cells.Item(87, 8).Value = 143407
$ws.Cells.Item(87, 10).Value = 139810.33
$ws.Cells.Item(87, 12).Value = 139810.33
$ws.Cells.Item(87, 14).Value = -142306.33
# Row 90
$ws.Cells.Item(90, 8).Value = 143407
$ws.Cells.Item(90, 10).Value = 139810.33
$ws.Cells.Item(90, 12).Value = 419430.99
$ws.Cells.Item(90, 14).Value = -431910.99
# Row 121
$ws.Cells.Item(121, 8).Value = 2451.2593
$ws.Cells.Item(121, 10).Value = 2514.7693
$ws.Cells.Item(121, 12).Value = 7544.3079
$ws.Cells.Item(121, 14).Value = -11038.3079
# Row 127
$ws.Cells.Item(127, 8).Value = 800.86664
$ws.Cells.Item(127, 9).Value = 643.8570999999999
$ws.Cells.Item(127, 10).Value = 2999
$ws.Cells.Item(127, 11).Value = 1931.5713
$ws.Cells.Item(127, 12).Value = 8997
$ws.Cells.Item(127, 13).Value = 3028.4287
$ws.Cells.Item(127, 14).Value = -18917
# Row 138
$ws.Cells.Item(138, 8).Value = 2491.87
$ws.Cells.Item(138, 9).Value = 1701.909
$ws.Cells.Item(138, 10).Value = 2880.9553
$ws.Cells.Item(138, 11).Value = 5105.727000000001
$ws.Cells.Item(138, 12).Value = 8642.865900000001
$ws.Cells.Item(138, 13).Value = 34.27299999999923
$ws.Cells.Item(138, 14).Value = -18922.8659

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Cells.Item(61, 8).Value = 2731.4194
$ws.Cells.Item(61, 9).Value = 2093
$ws.Cells.Item(61, 11).Value = 2093
$ws.Cells.Item(61, 13).Value = -1881
# Row 62
$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 14).ClearContents()
$ws.Cells.Item(62, 12).Value = 0
# Row 65
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 14).ClearContents()
$ws.Cells.Item(65, 12).Value = 0
# Row 88
$ws.Cells.Item(88, 8).Value = 2910.0833
$ws.Cells.Item(88, 9).Value = 2999.6667
$ws.Cells.Item(88, 10).Value = 2880.2222
$ws.Cells.Item(88, 11).Value = 2999.6667
$ws.Cells.Item(88, 12).Value = 2880.2222
$ws.Cells.Item(88, 13).Value = -2593.6667
$ws.Cells.Item(88, 14).Value = -3692.2222
# Row 91
$ws.Cells.Item(91, 8).Value = 2910.0833
$ws.Cells.Item(91, 9).Value = 2999.6667
$ws.Cells.Item(91, 10).Value = 2880.2222
$ws.Cells.Item(91, 11).Value = 2999.6667
$ws.Cells.Item(91, 12).Value = 2880.2222
$ws.Cells.Item(91, 13).Value = -1595.6667
$ws.Cells.Item(91, 14).Value = -5688.2222
# Row 102
$ws.Cells.Item(102, 8).Value = 1583.0526
$ws.Cells.Item(102, 9).Value = 1005.3333
$ws.Cells.Item(102, 10).Value = 3749.5
$ws.Cells.Item(102, 11).Value = 1005.3333
$ws.Cells.Item(102, 12).Value = 3749.5
$ws.Cells.Item(102, 13).Value = 616.6667
$ws.Cells.Item(102, 14).Value = -6993.5
# Row 132
$ws.Cells.Item(132, 8).Value = 180789.02
$ws.Cells.Item(132, 9).Value = 223948.42
$ws.Cells.Item(132, 11).Value = 671845.26
$ws.Cells.Item(132, 13).Value = -669315.26
# Row 136
$ws.Cells.Item(136, 8).Value = 2731.4194
$ws.Cells.Item(136, 9).Value = 2093
$ws.Cells.Item(136, 11).Value = 6279
$ws.Cells.Item(136, 13).Value = -3729

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Cells.Item(86, 8).Value = 3570.375
$ws.Cells.Item(86, 9).Value = 3162.2
$ws.Cells.Item(86, 10).Value = 4250.6665
$ws.Cells.Item(86, 11).Value = 3162.2
$ws.Cells.Item(86, 12).Value = 4250.6665
$ws.Cells.Item(86, 13).Value = -2039.2
$ws.Cells.Item(86, 14).Value = -6496.6665
# Row 89
$ws.Cells.Item(89, 8).Value = 3570.375
$ws.Cells.Item(89, 9).Value = 3162.2
$ws.Cells.Item(89, 10).Value = 4250.6665
$ws.Cells.Item(89, 11).Value = 15811
$ws.Cells.Item(89, 12).Value = 21253.3325
$ws.Cells.Item(89, 13).Value = -10195
$ws.Cells.Item(89, 14).Value = -32485.3325
# Row 99
$ws.Cells.Item(99, 8).Value = 2523.0435
$ws.Cells.Item(99, 9).Value = 2195.6667
$ws.Cells.Item(99, 10).Value = 3136.875
$ws.Cells.Item(99, 11).Value = 2195.6667
$ws.Cells.Item(99, 12).Value = 3136.875
$ws.Cells.Item(99, 13).Value = -697.6667000000002
$ws.Cells.Item(99, 14).Value = -6132.875
# Row 105
$ws.Cells.Item(105, 8).Value = 3047.2
$ws.Cells.Item(105, 9).Value = 2802.3333
$ws.Cells.Item(105, 11).Value = 2802.3333
$ws.Cells.Item(105, 13).Value = -1055.3333
# Row 134
$ws.Cells.Item(134, 8).Value = 2103787.8
$ws.Cells.Item(134, 9).Value = 2749185.5
$ws.Cells.Item(134, 10).Value = 6245.625
$ws.Cells.Item(134, 11).Value = 8247556.5
$ws.Cells.Item(134, 12).Value = 18736.875
$ws.Cells.Item(134, 13).Value = -8245021.5
$ws.Cells.Item(134, 14).Value = -23806.875

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 68
$ws.Cells.Item(68, 8).Value = 31500
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).ClearContents()
# Row 71
$ws.Cells.Item(71, 8).Value = 31500
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).ClearContents()
# Row 86
$ws.Cells.Item(86, 8).Value = 35209.418
$ws.Cells.Item(86, 9).Value = 28460.54
$ws.Cells.Item(86, 11).Value = 28460.54
$ws.Cells.Item(86, 13).Value = -27337.54
# Row 89
$ws.Cells.Item(89, 8).Value = 35209.418
$ws.Cells.Item(89, 9).Value = 28460.54
$ws.Cells.Item(89, 11).Value = 142302.7
$ws.Cells.Item(89, 13).Value = -136686.7
# Row 107
$ws.Cells.Item(107, 8).Value = 1178.1666
$ws.Cells.Item(107, 9).Value = 740
$ws.Cells.Item(107, 10).Value = 1528.7
$ws.Cells.Item(107, 11).Value = 740
$ws.Cells.Item(107, 12).Value = 1528.7
$ws.Cells.Item(107, 13).Value = 1180
$ws.Cells.Item(107, 14).Value = -5368.7
# Row 118
$ws.Cells.Item(118, 8).Value = 111490
$ws.Cells.Item(118, 10).Value = 103980
$ws.Cells.Item(118, 12).Value = 103980
$ws.Cells.Item(118, 14).Value = -107294

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 13).ClearContents()
# Row 68
$ws.Cells.Item(68, 8).Value = 1166.6
$ws.Cells.Item(68, 10).Value = 1166.6
$ws.Cells.Item(68, 12).Value = 3499.8
$ws.Cells.Item(68, 14).Value = -5121.799999999999
# Row 71
$ws.Cells.Item(71, 8).Value = 1166.6
$ws.Cells.Item(71, 10).Value = 1166.6
$ws.Cells.Item(71, 12).Value = 10499.4
$ws.Cells.Item(71, 14).Value = -18611.4
# Row 82
$ws.Cells.Item(82, 8).Value = 1500
$ws.Cells.Item(82, 9).Value = 1500
$ws.Cells.Item(82, 11).Value = 4500
$ws.Cells.Item(82, 13).Value = -4094
# Row 85
$ws.Cells.Item(85, 8).Value = 1500
$ws.Cells.Item(85, 9).Value = 1500
$ws.Cells.Item(85, 11).Value = 4500
$ws.Cells.Item(85, 13).Value = -3096
# Row 132
$ws.Cells.Item(132, 8).Value = 980.2778
$ws.Cells.Item(132, 9).Value = 1750
$ws.Cells.Item(132, 10).Value = 935
$ws.Cells.Item(132, 11).Value = 15750
$ws.Cells.Item(132, 12).Value = 8415
$ws.Cells.Item(132, 13).Value = -13220
$ws.Cells.Item(132, 14).Value = -13475

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 24
$ws.Cells.Item(24, 8).Value = 125043660
$ws.Cells.Item(24, 9).Value = 1000000000
$ws.Cells.Item(24, 10).Value = 49904
$ws.Cells.Item(24, 11).Value = 1000000000
$ws.Cells.Item(24, 12).Value = 49904
$ws.Cells.Item(24, 13).Value = -999999827
$ws.Cells.Item(24, 14).Value = -50250
# Row 80
$ws.Cells.Item(80, 8).Value = 2168.3125
$ws.Cells.Item(80, 9).Value = 1779.9166
$ws.Cells.Item(80, 11).Value = 1779.9166
$ws.Cells.Item(80, 13).Value = -781.9166
# Row 83
$ws.Cells.Item(83, 8).Value = 2168.3125
$ws.Cells.Item(83, 9).Value = 1779.9166
$ws.Cells.Item(83, 11).Value = 8899.583000000001
$ws.Cells.Item(83, 13).Value = -3907.583000000001
# Row 97
$ws.Cells.Item(97, 8).Value = 997.6
$ws.Cells.Item(97, 9).Value = 669.2
$ws.Cells.Item(97, 11).Value = 669.2
$ws.Cells.Item(97, 13).Value = -173.2
# Row 102
$ws.Cells.Item(102, 8).Value = 2472.4375
$ws.Cells.Item(102, 9).Value = 2472.4375
$ws.Cells.Item(102, 11).Value = 2472.4375
$ws.Cells.Item(102, 13).Value = -850.4375
# Row 120
$ws.Cells.Item(120, 8).Value = 88229.336
$ws.Cells.Item(120, 10).Value = 88229.336
$ws.Cells.Item(120, 12).Value = 88229.336
$ws.Cells.Item(120, 14).Value = -97905.336

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 4200.4585
$ws.Cells.Item(7, 9).Value = 3826.389
$ws.Cells.Item(7, 10).Value = 5322.6665
$ws.Cells.Item(7, 11).Value = 3826.389
$ws.Cells.Item(7, 12).Value = 5322.6665
$ws.Cells.Item(7, 13).Value = -3714.389
$ws.Cells.Item(7, 14).Value = -5546.6665
# Row 40
$ws.Cells.Item(40, 8).Value = 4265.364
$ws.Cells.Item(40, 9).Value = 1991
$ws.Cells.Item(40, 11).Value = 1991
$ws.Cells.Item(40, 13).Value = -1855
# Row 81
$ws.Cells.Item(81, 8).Value = 129989
$ws.Cells.Item(81, 10).Value = 129989
$ws.Cells.Item(81, 12).Value = 129989
$ws.Cells.Item(81, 14).Value = -131985
# Row 84
$ws.Cells.Item(84, 8).Value = 129989
$ws.Cells.Item(84, 10).Value = 129989
$ws.Cells.Item(84, 12).Value = 389967
$ws.Cells.Item(84, 14).Value = -399951
# Row 126
$ws.Cells.Item(126, 8).Value = 4200.4585
$ws.Cells.Item(126, 9).Value = 3826.389
$ws.Cells.Item(126, 10).Value = 5322.6665
$ws.Cells.Item(126, 11).Value = 11479.167
$ws.Cells.Item(126, 12).Value = 15967.9995
$ws.Cells.Item(126, 13).Value = -9009.167000000001
$ws.Cells.Item(126, 14).Value = -20907.9995
# Row 132
$ws.Cells.Item(132, 8).Value = 46987.566
$ws.Cells.Item(132, 9).Value = 53196.4
$ws.Cells.Item(132, 11).Value = 159589.2
$ws.Cells.Item(132, 13).Value = -157059.2
# Row 139
$ws.Cells.Item(139, 8).Value = 119999.5
$ws.Cells.Item(139, 10).Value = 119999.5
$ws.Cells.Item(139, 12).Value = 119999.5
$ws.Cells.Item(139, 14).Value = -130279.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 18
$ws.Cells.Item(18, 8).Value = 18752.5
$ws.Cells.Item(18, 9).Value = 17499.75
$ws.Cells.Item(18, 10).Value = 20005.25
$ws.Cells.Item(18, 11).Value = 17499.75
$ws.Cells.Item(18, 12).Value = 20005.25
$ws.Cells.Item(18, 13).Value = -17326.75
$ws.Cells.Item(18, 14).Value = -20351.25
# Row 64
$ws.Cells.Item(64, 8).Value = 99982.664
$ws.Cells.Item(64, 10).Value = 99982.664
$ws.Cells.Item(64, 12).Value = 99982.664
$ws.Cells.Item(64, 14).Value = -100478.664
# Row 67
$ws.Cells.Item(67, 8).Value = 99982.664
$ws.Cells.Item(67, 10).Value = 99982.664
$ws.Cells.Item(67, 12).Value = 99982.664
$ws.Cells.Item(67, 14).Value = -101698.664
# Row 126
$ws.Cells.Item(126, 8).Value = 4151
$ws.Cells.Item(126, 9).Value = 4136.706
$ws.Cells.Item(126, 10).Value = 4199.6
$ws.Cells.Item(126, 11).Value = 12410.118
$ws.Cells.Item(126, 12).Value = 12598.8
$ws.Cells.Item(126, 13).Value = -9940.118
$ws.Cells.Item(126, 14).Value = -17538.8
# Row 132
$ws.Cells.Item(132, 8).Value = 31378.428
$ws.Cells.Item(132, 9).Value = 35988.867
$ws.Cells.Item(132, 10).Value = 3715.8
$ws.Cells.Item(132, 11).Value = 107966.601
$ws.Cells.Item(132, 12).Value = 11147.4
$ws.Cells.Item(132, 13).Value = -105436.601
$ws.Cells.Item(132, 14).Value = -16207.4
